# Edit script: reformulated MIP procurement model output.
# The two result sheets swap roles/names (pet_gourmet becomes the first sheet,
# patas_pack becomes the second), and the solved values for packings P4, P5 and
# P6 change to reflect the updated formulation/new constraints.

$wb = $excel.ActiveWorkbook

# Grab sheets by their current (pre-edit) tab position before any renaming.
$ws1 = $wb.Worksheets.Item(1)   # currently "patas_pack" -> becomes "pet_gourmet"
$ws2 = $wb.Worksheets.Item(2)   # currently "pet_gourmet" -> becomes "patas_pack"

# Swap the sheet names (via a temporary name to avoid a name collision).
$ws1.Name = "__tmp_swap__"
$ws2.Name = "patas_pack"
$ws1.Name = "pet_gourmet"

# --- Update values on the sheet now named "pet_gourmet" (tab 1) ---
$ws1.Cells.Item(1,4).Value2 = "Demand"
$ws1.Cells.Item(1,5).Value2 = "Transferred Quantity"
$ws1.Cells.Item(2,3).Value2 = 200
$ws1.Cells.Item(2,4).Value2 = 2000
$ws1.Cells.Item(2,5).Value2 = 1810
$ws1.Cells.Item(2,6).Value2 = 10
$ws1.Cells.Item(3,3).Value2 = 10
$ws1.Cells.Item(3,6).Value2 = 10
$ws1.Cells.Item(4,3).Value2 = 10
$ws1.Cells.Item(4,6).Value2 = 10
$ws1.Cells.Item(5,3).Value2 = 50
$ws1.Cells.Item(5,4).Value2 = 2000
$ws1.Cells.Item(5,5).Value2 = 1960
$ws1.Cells.Item(5,6).Value2 = 10
$ws1.Cells.Item(6,3).Value2 = 10
$ws1.Cells.Item(6,6).Value2 = 10
$ws1.Cells.Item(7,3).Value2 = 10
$ws1.Cells.Item(7,6).Value2 = 10
$ws1.Cells.Item(8,3).Value2 = 100
$ws1.Cells.Item(8,4).Value2 = 2500
$ws1.Cells.Item(8,5).Value2 = 2410
$ws1.Cells.Item(8,6).Value2 = 10
$ws1.Cells.Item(9,3).Value2 = 10
$ws1.Cells.Item(9,6).Value2 = 10
$ws1.Cells.Item(10,3).Value2 = 10
$ws1.Cells.Item(10,6).Value2 = 10
$ws1.Cells.Item(11,3).Value2 = 200
$ws1.Cells.Item(11,4).Value2 = 1000
$ws1.Cells.Item(11,5).Value2 = 2310
$ws1.Cells.Item(11,6).Value2 = 1510
$ws1.Cells.Item(12,3).Value2 = 1510
$ws1.Cells.Item(12,5).Value2 = 0
$ws1.Cells.Item(12,6).Value2 = 10
$ws1.Cells.Item(13,3).Value2 = 10
$ws1.Cells.Item(13,6).Value2 = 10
$ws1.Cells.Item(14,3).Value2 = 200
$ws1.Cells.Item(14,4).Value2 = 1000
$ws1.Cells.Item(14,5).Value2 = 810
$ws1.Cells.Item(14,6).Value2 = 10
$ws1.Cells.Item(15,3).Value2 = 10
$ws1.Cells.Item(15,5).Value2 = 2000
$ws1.Cells.Item(15,6).Value2 = 1010
$ws1.Cells.Item(16,3).Value2 = 1010
$ws1.Cells.Item(16,5).Value2 = 0
$ws1.Cells.Item(16,6).Value2 = 10
$ws1.Cells.Item(17,3).Value2 = 2900
$ws1.Cells.Item(17,4).Value2 = 2000
$ws1.Cells.Item(17,5).Value2 = 0
$ws1.Cells.Item(17,6).Value2 = 900
$ws1.Cells.Item(18,3).Value2 = 900
$ws1.Cells.Item(18,5).Value2 = 1110
$ws1.Cells.Item(18,6).Value2 = 10
$ws1.Cells.Item(19,3).Value2 = 10
$ws1.Cells.Item(19,6).Value2 = 10

# --- Update values on the sheet now named "patas_pack" (tab 2) ---
$ws2.Cells.Item(1,4).Value2 = "Transferred Quantity"
$ws2.Cells.Item(1,5).Value2 = "Acquired Quantity"
$ws2.Cells.Item(2,3).Value2 = 100
$ws2.Cells.Item(2,4).Value2 = 1810
$ws2.Cells.Item(2,5).Value2 = 1710
$ws2.Cells.Item(2,6).Value2 = 0
$ws2.Cells.Item(3,3).Value2 = 0
$ws2.Cells.Item(3,6).Value2 = 0
$ws2.Cells.Item(4,3).Value2 = 0
$ws2.Cells.Item(4,6).Value2 = 0
$ws2.Cells.Item(5,3).Value2 = 100
$ws2.Cells.Item(5,4).Value2 = 1960
$ws2.Cells.Item(5,5).Value2 = 1860
$ws2.Cells.Item(5,6).Value2 = 0
$ws2.Cells.Item(6,3).Value2 = 0
$ws2.Cells.Item(6,6).Value2 = 0
$ws2.Cells.Item(7,3).Value2 = 0
$ws2.Cells.Item(7,6).Value2 = 0
$ws2.Cells.Item(8,3).Value2 = 50
$ws2.Cells.Item(8,4).Value2 = 2410
$ws2.Cells.Item(8,5).Value2 = 2360
$ws2.Cells.Item(8,6).Value2 = 0
$ws2.Cells.Item(9,3).Value2 = 0
$ws2.Cells.Item(9,6).Value2 = 0
$ws2.Cells.Item(10,3).Value2 = 0
$ws2.Cells.Item(10,6).Value2 = 0
$ws2.Cells.Item(11,3).Value2 = 100
$ws2.Cells.Item(11,4).Value2 = 2310
$ws2.Cells.Item(11,5).Value2 = 2210
$ws2.Cells.Item(11,6).Value2 = 0
$ws2.Cells.Item(12,3).Value2 = 0
$ws2.Cells.Item(12,4).Value2 = 0
$ws2.Cells.Item(12,5).Value2 = 0
$ws2.Cells.Item(12,6).Value2 = 0
$ws2.Cells.Item(13,3).Value2 = 0
$ws2.Cells.Item(13,6).Value2 = 0
$ws2.Cells.Item(14,3).Value2 = 100
$ws2.Cells.Item(14,4).Value2 = 810
$ws2.Cells.Item(14,5).Value2 = 710
$ws2.Cells.Item(14,6).Value2 = 0
$ws2.Cells.Item(15,3).Value2 = 0
$ws2.Cells.Item(15,4).Value2 = 2000
$ws2.Cells.Item(15,5).Value2 = 2000
$ws2.Cells.Item(15,6).Value2 = 0
$ws2.Cells.Item(16,3).Value2 = 0
$ws2.Cells.Item(16,4).Value2 = 0
$ws2.Cells.Item(16,5).Value2 = 0
$ws2.Cells.Item(16,6).Value2 = 0
$ws2.Cells.Item(17,4).Value2 = 0
$ws2.Cells.Item(17,5).Value2 = 0
$ws2.Cells.Item(17,6).Value2 = 50
$ws2.Cells.Item(18,3).Value2 = 50
$ws2.Cells.Item(18,4).Value2 = 1110
$ws2.Cells.Item(18,5).Value2 = 1060
$ws2.Cells.Item(18,6).Value2 = 0
$ws2.Cells.Item(19,3).Value2 = 0
$ws2.Cells.Item(19,6).Value2 = 0
